# Lta-Ltbr.xlsx: refresh NATMI TPM-derived LR-pair metrics (new TPM run).
# - Numeric columns E:T for data rows 2-21 get new computed values.
# - Column A (Sending cluster) for rows 17-21 is relabeled from "MuSCs" to
#   "Resolving-Mac" (cluster identity correction that came with the re-run).
# - Column D (Target cluster) text is unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 20,16
$arr[0,0] = 1
$arr[0,1] = 0.3333333333333333
$arr[0,2] = 0.07057933333333334
$arr[0,3] = 0.211738
$arr[0,4] = 0.1281663575568867
$arr[0,5] = 0.1281663575568867
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 9.071155000000001
$arr[0,9] = 27.213465
$arr[0,10] = 0.1321733179750106
$arr[0,11] = 0.135675019324872
$arr[0,12] = 0.6402360724633335
$arr[0,13] = 5.762124652170001
$arr[0,14] = 0.01694017273106529
$arr[0,15] = 0.01738897303832905
$arr[1,0] = 1
$arr[1,1] = 0.3333333333333333
$arr[1,2] = 0.07057933333333334
$arr[1,3] = 0.211738
$arr[1,4] = 0.1281663575568867
$arr[1,5] = 0.1281663575568867
$arr[1,6] = 3
$arr[1,7] = 1
$arr[1,8] = 20.60908733333333
$arr[1,9] = 61.82726199999999
$arr[1,10] = 0.3002893736556623
$arr[1,11] = 0.3082450164524775
$arr[1,12] = 1.454575644595111
$arr[1,13] = 13.091180801356
$arr[1,14] = 0.03848699523448518
$arr[1,15] = 0.03950664099377665
$arr[2,0] = 1
$arr[2,1] = 0.3333333333333333
$arr[2,2] = 0.07057933333333334
$arr[2,3] = 0.211738
$arr[2,4] = 0.1281663575568867
$arr[2,5] = 0.1281663575568867
$arr[2,6] = 3
$arr[2,7] = 1
$arr[2,8] = 17.81090666666666
$arr[2,9] = 53.43272
$arr[2,10] = 0.2595178486396241
$arr[2,11] = 0.2663933210482557
$arr[2,12] = 1.257081918595556
$arr[2,13] = 11.31373726736
$arr[2,14] = 0.03326145738114007
$arr[2,15] = 0.03414266163623726
$arr[3,0] = 1
$arr[3,1] = 0.3333333333333333
$arr[3,2] = 0.07057933333333334
$arr[3,3] = 0.211738
$arr[3,4] = 0.1281663575568867
$arr[3,5] = 0.1281663575568867
$arr[3,6] = 2
$arr[3,7] = 1
$arr[3,8] = 5.313972
$arr[3,9] = 10.627944
$arr[3,10] = 0.07742843230727542
$arr[3,11] = 0.05298650897942091
$arr[3,12] = 0.375056601112
$arr[3,13] = 2.250339606672
$arr[3,14] = 0.00992372014016346
$arr[3,15] = 0.00679108785554765
$arr[4,0] = 1
$arr[4,1] = 0.3333333333333333
$arr[4,2] = 0.07057933333333334
$arr[4,3] = 0.211738
$arr[4,4] = 0.1281663575568867
$arr[4,5] = 0.1281663575568867
$arr[4,6] = 3
$arr[4,7] = 1
$arr[4,8] = 15.825637
$arr[4,9] = 47.476911
$arr[4,10] = 0.2305910274224278
$arr[4,11] = 0.2367001341949739
$arr[4,12] = 1.116962909035333
$arr[4,13] = 10.052666181318
$arr[4,14] = 0.02955401207003274
$arr[4,15] = 0.03033699403299609
$arr[5,0] = 1
$arr[5,1] = 0.3333333333333333
$arr[5,2] = 0.0279
$arr[5,3] = 0.0837
$arr[5,4] = 0.05066414213561767
$arr[5,5] = 0.05066414213561768
$arr[5,6] = 3
$arr[5,7] = 1
$arr[5,8] = 9.071155000000001
$arr[5,9] = 27.213465
$arr[5,10] = 0.1321733179750106
$arr[5,11] = 0.135675019324872
$arr[5,12] = 0.2530852245
$arr[5,13] = 2.2777670205
$arr[5,14] = 0.006696447768422127
$arr[5,15] = 0.006873858463327989
$arr[6,0] = 1
$arr[6,1] = 0.3333333333333333
$arr[6,2] = 0.0279
$arr[6,3] = 0.0837
$arr[6,4] = 0.05066414213561767
$arr[6,5] = 0.05066414213561768
$arr[6,6] = 3
$arr[6,7] = 1
$arr[6,8] = 20.60908733333333
$arr[6,9] = 61.82726199999999
$arr[6,10] = 0.3002893736556623
$arr[6,11] = 0.3082450164524775
$arr[6,12] = 0.5749935365999999
$arr[6,13] = 5.174941829399999
$arr[6,14] = 0.01521390350870608
$arr[6,15] = 0.01561696932614413
$arr[7,0] = 1
$arr[7,1] = 0.3333333333333333
$arr[7,2] = 0.0279
$arr[7,3] = 0.0837
$arr[7,4] = 0.05066414213561767
$arr[7,5] = 0.05066414213561768
$arr[7,6] = 3
$arr[7,7] = 1
$arr[7,8] = 17.81090666666666
$arr[7,9] = 53.43272
$arr[7,10] = 0.2595178486396241
$arr[7,11] = 0.2663933210482557
$arr[7,12] = 0.4969242959999999
$arr[7,13] = 4.472318663999999
$arr[7,14] = 0.01314824917020763
$arr[7,15] = 0.01349658908156806
$arr[8,0] = 1
$arr[8,1] = 0.3333333333333333
$arr[8,2] = 0.0279
$arr[8,3] = 0.0837
$arr[8,4] = 0.05066414213561767
$arr[8,5] = 0.05066414213561768
$arr[8,6] = 2
$arr[8,7] = 1
$arr[8,8] = 5.313972
$arr[8,9] = 10.627944
$arr[8,10] = 0.07742843230727542
$arr[8,11] = 0.05298650897942091
$arr[8,12] = 0.1482598188
$arr[8,13] = 0.8895589127999999
$arr[8,14] = 0.003922845099753853
$arr[8,15] = 0.002684516022203564
$arr[9,0] = 1
$arr[9,1] = 0.3333333333333333
$arr[9,2] = 0.0279
$arr[9,3] = 0.0837
$arr[9,4] = 0.05066414213561767
$arr[9,5] = 0.05066414213561768
$arr[9,6] = 3
$arr[9,7] = 1
$arr[9,8] = 15.825637
$arr[9,9] = 47.476911
$arr[9,10] = 0.2305910274224278
$arr[9,11] = 0.2367001341949739
$arr[9,12] = 0.4415352722999999
$arr[9,13] = 3.9738174507
$arr[9,14] = 0.01168269658852799
$arr[9,15] = 0.01199220924237394
$arr[10,0] = 1
$arr[10,1] = 0.3333333333333333
$arr[10,2] = 0.128589
$arr[10,3] = 0.385767
$arr[10,4] = 0.2335072176730087
$arr[10,5] = 0.2335072176730087
$arr[10,6] = 3
$arr[10,7] = 1
$arr[10,8] = 9.071155000000001
$arr[10,9] = 27.213465
$arr[10,10] = 0.1321733179750106
$arr[10,11] = 0.135675019324872
$arr[10,12] = 1.166450750295
$arr[10,13] = 10.498056752655
$arr[10,14] = 0.03086342373095459
$arr[10,15] = 0.03168109627028254
$arr[11,0] = 1
$arr[11,1] = 0.3333333333333333
$arr[11,2] = 0.128589
$arr[11,3] = 0.385767
$arr[11,4] = 0.2335072176730087
$arr[11,5] = 0.2335072176730087
$arr[11,6] = 3
$arr[11,7] = 1
$arr[11,8] = 20.60908733333333
$arr[11,9] = 61.82726199999999
$arr[11,10] = 0.3002893736556623
$arr[11,11] = 0.3082450164524775
$arr[11,12] = 2.650101931106
$arr[11,13] = 23.850917379954
$arr[11,14] = 0.07011973613910417
$arr[11,15] = 0.07197743615338881
$arr[12,0] = 1
$arr[12,1] = 0.3333333333333333
$arr[12,2] = 0.128589
$arr[12,3] = 0.385767
$arr[12,4] = 0.2335072176730087
$arr[12,5] = 0.2335072176730087
$arr[12,6] = 3
$arr[12,7] = 1
$arr[12,8] = 17.81090666666666
$arr[12,9] = 53.43272
$arr[12,10] = 0.2595178486396241
$arr[12,11] = 0.2663933210482557
$arr[12,12] = 2.29028667736
$arr[12,13] = 20.61258009624
$arr[12,14] = 0.06059929077232361
$arr[12,15] = 0.06220476320465074
$arr[13,0] = 1
$arr[13,1] = 0.3333333333333333
$arr[13,2] = 0.128589
$arr[13,3] = 0.385767
$arr[13,4] = 0.2335072176730087
$arr[13,5] = 0.2335072176730087
$arr[13,6] = 2
$arr[13,7] = 1
$arr[13,8] = 5.313972
$arr[13,9] = 10.627944
$arr[13,10] = 0.07742843230727542
$arr[13,11] = 0.05298650897942091
$arr[13,12] = 0.683318345508
$arr[13,13] = 4.099910073048
$arr[13,14] = 0.01808009779685478
$arr[13,15] = 0.01237273228599047
$arr[14,0] = 1
$arr[14,1] = 0.3333333333333333
$arr[14,2] = 0.128589
$arr[14,3] = 0.385767
$arr[14,4] = 0.2335072176730087
$arr[14,5] = 0.2335072176730087
$arr[14,6] = 3
$arr[14,7] = 1
$arr[14,8] = 15.825637
$arr[14,9] = 47.476911
$arr[14,10] = 0.2305910274224278
$arr[14,11] = 0.2367001341949739
$arr[14,12] = 2.035002836193
$arr[14,13] = 18.315025525737
$arr[14,14] = 0.05384466923377155
$arr[14,15] = 0.05527118975869615
$arr[15,0] = 3
$arr[15,1] = 1
$arr[15,2] = 0.323617
$arr[15,3] = 0.970851
$arr[15,4] = 0.5876622826344869
$arr[15,5] = 0.5876622826344869
$arr[15,6] = 3
$arr[15,7] = 1
$arr[15,8] = 9.071155000000001
$arr[15,9] = 27.213465
$arr[15,10] = 0.1321733179750106
$arr[15,11] = 0.135675019324872
$arr[15,12] = 2.935579967635
$arr[15,13] = 26.420219708715
$arr[15,14] = 0.0776732737445686
$arr[15,15] = 0.07973109155293238
$arr[16,0] = 3
$arr[16,1] = 1
$arr[16,2] = 0.323617
$arr[16,3] = 0.970851
$arr[16,4] = 0.5876622826344869
$arr[16,5] = 0.5876622826344869
$arr[16,6] = 3
$arr[16,7] = 1
$arr[16,8] = 20.60908733333333
$arr[16,9] = 61.82726199999999
$arr[16,10] = 0.3002893736556623
$arr[16,11] = 0.3082450164524775
$arr[16,12] = 6.669451015551332
$arr[16,13] = 60.02505913996199
$arr[16,14] = 0.1764687387733669
$arr[16,15] = 0.1811439699791679
$arr[17,0] = 3
$arr[17,1] = 1
$arr[17,2] = 0.323617
$arr[17,3] = 0.970851
$arr[17,4] = 0.5876622826344869
$arr[17,5] = 0.5876622826344869
$arr[17,6] = 3
$arr[17,7] = 1
$arr[17,8] = 17.81090666666666
$arr[17,9] = 53.43272
$arr[17,10] = 0.2595178486396241
$arr[17,11] = 0.2663933210482557
$arr[17,12] = 5.763912182746665
$arr[17,13] = 51.87520964471999
$arr[17,14] = 0.1525088513159528
$arr[17,15] = 0.1565493071257997
$arr[18,0] = 3
$arr[18,1] = 1
$arr[18,2] = 0.323617
$arr[18,3] = 0.970851
$arr[18,4] = 0.5876622826344869
$arr[18,5] = 0.5876622826344869
$arr[18,6] = 2
$arr[18,7] = 1
$arr[18,8] = 5.313972
$arr[18,9] = 10.627944
$arr[18,10] = 0.07742843230727542
$arr[18,11] = 0.05298650897942091
$arr[18,12] = 1.719691676724
$arr[18,13] = 10.318150060344
$arr[18,14] = 0.04550176927050333
$arr[18,15] = 0.03113817281567923
$arr[19,0] = 3
$arr[19,1] = 1
$arr[19,2] = 0.323617
$arr[19,3] = 0.970851
$arr[19,4] = 0.5876622826344869
$arr[19,5] = 0.5876622826344869
$arr[19,6] = 3
$arr[19,7] = 1
$arr[19,8] = 15.825637
$arr[19,9] = 47.476911
$arr[19,10] = 0.2305910274224278
$arr[19,11] = 0.2367001341949739
$arr[19,12] = 5.121445169029
$arr[19,13] = 46.093006521261
$arr[19,14] = 0.1355096495300955
$arr[19,15] = 0.1390997411609078

$ws.Range("E2:T21").Value = $arr

$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("A18").Value = "Resolving-Mac"
$ws.Range("A19").Value = "Resolving-Mac"
$ws.Range("A20").Value = "Resolving-Mac"
$ws.Range("A21").Value = "Resolving-Mac"

Write-Output "done"